$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A15").Value = "#100013"
$ws.Range("B15").Value = "Pom cleaning and complete README.ME file"

$ws.Range("D27").Select()
